$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
$ws.Hyperlinks.Add($ws.Range("B20"), "mailto:test@test.com") | Out-Null
